# Updated cryptos list on Wed Sep  6 22:17:05 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    # Force the cell to stay a text/string value (not auto-converted to a
    # number by Excel's type inference), then drop back to the default
    # "Normal" style so no stray number-format style sticks to the cell.
    $r = $ws.Range($rangeAddress)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "25.782.63"
Set-TextValue "E2" "  -0.21%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.635.24"
Set-TextValue "E3" "  +0.01%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.12%  "

# Row 5 - BNB
Set-TextValue "D5" "215.18"
Set-TextValue "E5" "  +0.17%  "

# Row 6 - XRP
Set-TextValue "E6" "  -0.56%  "

# Row 7 - USDC
Set-TextValue "E7" "  -0.14%  "

# Row 8 - Cardano
Set-TextValue "E8" "  +0.34%  "

# Row 10 - Solana
Set-TextValue "D10" "19.64"
Set-TextValue "E10" "  -3.84%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0789"
Set-TextValue "E11" "  +1.23%  "

# Row 12 - Polkadot
Set-TextValue "E12" "  -0.04%  "

# Row 13 - now WrappedliquidstakedEther2.0 (was WrappedEther)
Set-TextValue "B13" "WrappedliquidstakedEther2.0"
Set-TextValue "C13" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D13" "1.860.28"
Set-TextValue "E13" "  -0.05%  "

# Row 14 - now WrappedEther (was WrappedliquidstakedEther2.0)
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.620.61"
Set-TextValue "E14" "  -1.21%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.559"
Set-TextValue "E15" "  -0.15%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.0₃0767"
Set-TextValue "E16" "  +0.06%  "

# Row 17 - Litecoin
Set-TextValue "D17" "62.88"
Set-TextValue "E17" "  -0.54%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "25.788.78"
Set-TextValue "E18" "  -0.23%  "

# Row 19 - Dai
Set-TextValue "E19" "  -0.20%  "

# Row 20 - Uniswap
Set-TextValue "E20" "  +1.62%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "194.57"
Set-TextValue "E21" "  +0.81%  "

# Row 22 - Avalanche
Set-TextValue "D22" "9.97"
Set-TextValue "E22" "  +0.10%  "

# Row 23 - Chainlink
Set-TextValue "E23" "  +2.32%  "

# Row 24 - BinanceUSD
Set-TextValue "E24" "  -0.16%  "

# Row 25 - Toncoin
Set-TextValue "D25" "1.83"
Set-TextValue "E25" "  +3.77%  "

# Row 26 - Monero
Set-TextValue "D26" "142.86"
Set-TextValue "E26" "  +3.49%  "

# Row 27 - Stellar
Set-TextValue "E27" "  -0.06%  "

# Row 28 - Cosmos
Set-TextValue "D28" "6.90"
Set-TextValue "E28" "  +0.78%  "

# Row 29 - EthereumClassic
Set-TextValue "E29" "  +0.42%  "

# Row 30 - PancakeSwap
Set-TextValue "E30" "  -0.13%  "

# Row 31 - Hedera
Set-TextValue "E31" "  -0.06%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "D32" "3.34"
Set-TextValue "E32" "  +1.05%  "

# Row 33 - Filecoin
Set-TextValue "E33" "  +0.24%  "

# Row 34 - LidoDAOToken
Set-TextValue "D34" "1.58"
Set-TextValue "E34" "  +0.98%  "

# Row 35 - HuobiToken
Set-TextValue "E35" "  -0.03%  "

# Row 36 - ARBITRUM
Set-TextValue "D36" "0.905"
Set-TextValue "E36" "  +0.30%  "

# Row 37 - Maker
Set-TextValue "D37" "1.129.28"
Set-TextValue "E37" "  -0.30%  "

# Row 38 - MXToken
Set-TextValue "E38" "  -1.90%  "

# Row 39 - ImmutableX
Set-TextValue "E39" "  -1.76%  "

# Row 40 - VeChain
Set-TextValue "E40" "  -0.20%  "

# Row 41 - PaxDollar
Set-TextValue "E41" "  +0.51%  "

# Row 42 - FraxShare
Set-TextValue "D42" "5.59"
Set-TextValue "E42" "  +2.34%  "

# Row 43 - Quant
Set-TextValue "E43" "  +1.04%  "

# Row 44 - TrustWalletToken
Set-TextValue "D44" "0.807"
Set-TextValue "E44" "  +0.86%  "

# Row 45 - RocketPoolETH
Set-TextValue "D45" "1.769.76"
Set-TextValue "E45" "  -0.23%  "

# Row 46 - BabyDogeCoin
Set-TextValue "D46" "0.0₆0110"
Set-TextValue "E46" "  -0.32%  "

# Row 47 - Aave
Set-TextValue "D47" "55.24"
Set-TextValue "E47" "  -0.56%  "

# Row 48 - Mantle
Set-TextValue "E48" "  -2.33%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.0503"
Set-TextValue "E49" "  -0.04%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "7.56"
Set-TextValue "E50" "  -2.90%  "

# Row 51 - now Frax (was SynthetixNetwork)
Set-TextValue "B51" "Frax"
Set-TextValue "C51" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D51" "1.00"
Set-TextValue "E51" "  +0.08%  "
